# Final commit: update the IMG attachment URL and refresh the active
# cell selection/scroll position on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Discord CDN image URL stored in L2 (shared string).
$ws.Range("L2").Value = "https://cdn.discordapp.com/attachments/1162451241872412901/1169225366959837214/IMG-20231101-WA0018.jpg?ex=6554a0fa&is=65422bfa&hm=220a0402abcc6ca88f3ef35be5638ff9b15a071285d00af17bfde03f88461d07&"

# Move the selection to G9 (this also resets any custom topLeftCell
# scroll position back to the natural view).
$ws.Range("G9").Select()
